# Move the "Message Format: ToString() vs Dump" slide (originally slide 14)
# to slide position 12, pushing the two picture slides that were at
# positions 12-13 down to positions 13-14.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(14)
$s.MoveTo(12)
